$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.083.46"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.638.31"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D5").Value = "'214.07"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "'0.5257"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "'0.06315"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("D11").Value = "'0.07666"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "1.629.09"
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").Value = "'4.427"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "1.860.61"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "'0.5511"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "0.0₅8148"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'65.17"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "26.078.82"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'4.692"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "'188.63"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "'6.169"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'146.30"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "'0.1219"
$ws.Range("D27").Value = "'7.409"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'1.406"
$ws.Range("E29").Value = "  +3.70%  "
$ws.Range("D30").Value = "'0.06002"
$ws.Range("E30").Value = "  -4.38%  "
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "'3.411"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "'1.642"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "'0.9877"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'2.760"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'0.5735"
$ws.Range("E38").Value = "  -5.31%  "
$ws.Range("D39").Value = "'0.01617"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "'0.8555"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "1.041.39"
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "'5.698"
$ws.Range("E43").Value = "  -7.00%  "
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "1.787.71"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("D47").Value = "'55.49"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'0.9980"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "'8.060"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'0.05173"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "'0.4223"
$ws.Range("E51").Value = "  -0.56%  "
